$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("O2").Value = 1.17
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 1.57
$ws.Range("R2").Value = 2.35

# Row 3 updates
$ws.Range("G3").Value = 2.05
$ws.Range("I3").Value = 3.9
$ws.Range("J3").Value = 2.88
$ws.Range("M3").Value = 1.13
$ws.Range("N3").Value = 6
$ws.Range("AN3").Value = 3.75
$ws.Range("AO3").Value = 12
$ws.Range("AZ3").Value = 101

# Row 4 updates
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 13
$ws.Range("O4").Value = 1.25
$ws.Range("P4").Value = 3.75
$ws.Range("Q4").Value = 1.88
$ws.Range("R4").Value = 1.98
